$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price data that must stay text (matches the source formatting,
# e.g. "98.588.21" / "255.61"); force a text number-format on the cells whose new
# value would otherwise be auto-detected as a number by the Value setter.
$textCells = @("D5", "D6", "D7", "D8", "D10", "D13", "D15", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D26", "D27", "D28", "D31", "D32", "D33", "D34", "D35", "D36", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '98.588.21'
$ws.Range('E2').Value = '  +1.49%  '
$ws.Range('D3').Value = '3.306.26'
$ws.Range('E3').Value = '  -0.36%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '255.61'
$ws.Range('E5').Value = '  +2.40%  '
$ws.Range('D6').Value = '624.60'
$ws.Range('E6').Value = '  +0.41%  '
$ws.Range('D7').Value = '1.46'
$ws.Range('E7').Value = '  +31.75%  '
$ws.Range('D8').Value = '0.401'
$ws.Range('E8').Value = '  +4.24%  '
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('D10').Value = '0.949'
$ws.Range('E10').Value = '  +20.48%  '
$ws.Range('D11').Value = '3.304.24'
$ws.Range('E11').Value = '  -0.46%  '
$ws.Range('E12').Value = '  +0.49%  '
$ws.Range('D13').Value = '39.05'
$ws.Range('E13').Value = '  +10.35%  '
$ws.Range('D14').Value = '98.262.86'
$ws.Range('E14').Value = '  +1.57%  '
$ws.Range('D15').Value = '0.0000249'
$ws.Range('E15').Value = '  +1.41%  '
$ws.Range('D16').Value = '3.922.27'
$ws.Range('E16').Value = '  +0.15%  '
$ws.Range('D17').Value = '5.49'
$ws.Range('E17').Value = '  -0.50%  '
$ws.Range('D18').Value = '3.302.31'
$ws.Range('E18').Value = '  +0.24%  '
$ws.Range('D19').Value = '3.48'
$ws.Range('E19').Value = '  -1.87%  '
$ws.Range('D20').Value = '15.48'
$ws.Range('E20').Value = '  +3.37%  '
$ws.Range('D21').Value = '6.37'
$ws.Range('E21').Value = '  +9.05%  '
$ws.Range('D22').Value = '484.44'
$ws.Range('E22').Value = '  +0.48%  '
$ws.Range('D23').Value = '9.45'
$ws.Range('E23').Value = '  +1.85%  '
$ws.Range('D24').Value = '0.0000203'
$ws.Range('E24').Value = '  -1.77%  '
$ws.Range('E25').Value = '  -0.67%  '
$ws.Range('D26').Value = '88.69'
$ws.Range('E26').Value = '  +1.00%  '
$ws.Range('D27').Value = '12.09'
$ws.Range('E27').Value = '  -0.49%  '
$ws.Range('B28').Value = 'Stellar'
$ws.Range('C28').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D28').Value = '0.302'
$ws.Range('E28').Value = '  +25.37%  '
$ws.Range('B29').Value = 'WrappedeETH'
$ws.Range('C29').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D29').Value = '3.478.55'
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('E30').Value = '  -0.18%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = '0.137'
$ws.Range('E31').Value = '  +11.97%  '
$ws.Range('B32').Value = 'Cronos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D32').Value = '0.188'
$ws.Range('E32').Value = '  +2.87%  '
$ws.Range('D33').Value = '10.11'
$ws.Range('E33').Value = '  +9.51%  '
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('D35').Value = '27.98'
$ws.Range('E35').Value = '  +2.54%  '
$ws.Range('D36').Value = '7.18'
$ws.Range('E36').Value = '  -3.15%  '
$ws.Range('E37').Value = '  -2.32%  '
$ws.Range('E38').Value = '  +0.81%  '
$ws.Range('E39').Value = '  +2.87%  '
$ws.Range('D40').Value = '24.82'
$ws.Range('E40').Value = '  +0.20%  '
$ws.Range('D41').Value = '490.59'
$ws.Range('E41').Value = '  -3.12%  '
$ws.Range('B42').Value = 'MantraDAO'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D42').Value = '3.69'
$ws.Range('E42').Value = '  +5.34%  '
$ws.Range('B43').Value = 'Fetch.AI'
$ws.Range('C43').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D43').Value = '1.23'
$ws.Range('E43').Value = '  -4.01%  '
$ws.Range('D44').Value = '0.795'
$ws.Range('E44').Value = '  -1.20%  '
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('D46').Value = '3.13'
$ws.Range('E46').Value = '  -4.10%  '
$ws.Range('B47').Value = 'Monero'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D47').Value = '158.25'
$ws.Range('E47').Value = '  -1.67%  '
$ws.Range('D48').Value = '1.93'
$ws.Range('E48').Value = '  +0.83%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').Value = '7.35'
$ws.Range('E49').Value = '  +16.00%  '
$ws.Range('B50').Value = 'Filecoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D50').Value = '4.74'
$ws.Range('E50').Value = '  +4.63%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = '0.848'
$ws.Range('E51').Value = '  +7.08%  '
